# Generate Report for Handoff
# Updates the "b.md" row across the Overview, zh-cn and de-de sheets to
# reflect that the file is now ready for handoff, with a fresh handoff
# file name / timestamp for each locale.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# ----------------------------------------------------------------------
# Overview sheet: row 3 corresponds to b.md
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsOverview.Range("D3").Value = "2016-03-24 04:34:43"

# ----------------------------------------------------------------------
# zh-cn sheet: row 3 corresponds to b.md
# ----------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$zhCnHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"

$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("D3").Value = $zhCnHandoffFile
$wsZhCn.Range("E3").Value = "2016-03-24 04:34:39"

foreach ($hl in $wsZhCn.Hyperlinks) {
    if ($hl.Range.Address() -eq '$D$3') {
        $hl.TextToDisplay = $zhCnHandoffFile
    }
}

# ----------------------------------------------------------------------
# de-de sheet: row 3 corresponds to b.md
# ----------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$deDeHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"

$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("D3").Value = $deDeHandoffFile
$wsDeDe.Range("E3").Value = "2016-03-24 04:34:43"

foreach ($hl in $wsDeDe.Hyperlinks) {
    if ($hl.Range.Address() -eq '$D$3') {
        $hl.TextToDisplay = $deDeHandoffFile
    }
}
